$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.401.52"
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").Value = "3.579.96"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "601.15"
$ws.Range("D6").Value = "135.42"
$ws.Range("E6").Value = "  -3.66%  "
$ws.Range("D7").Value = "3.579.37"
$ws.Range("E7").Value = "  +0.44%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "0.495"
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("E10").Value = "  -1.41%  "
$ws.Range("D11").Value = "7.18"
$ws.Range("E11").Value = "  +1.74%  "
$ws.Range("D12").Value = "0.391"
$ws.Range("E12").Value = "  -1.47%  "
$ws.Range("D13").Value = "4.183.09"
$ws.Range("E13").Value = "  +0.29%  "
$ws.Range("D14").Value = "0.0000185"
$ws.Range("E14").Value = "  -1.88%  "
$ws.Range("D15").Value = "27.51"
$ws.Range("E15").Value = "  +1.14%  "
$ws.Range("D16").Value = "3.573.10"
$ws.Range("E16").Value = "  +0.25%  "
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("D18").Value = "65.410.61"
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("D19").Value = "10.08"
$ws.Range("E19").Value = "  -2.66%  "
$ws.Range("D20").Value = "14.50"
$ws.Range("E20").Value = "  +1.49%  "
$ws.Range("D21").Value = "5.88"
$ws.Range("E21").Value = "  -0.54%  "
$ws.Range("D22").Value = "393.38"
$ws.Range("E22").Value = "  -1.05%  "
$ws.Range("D23").Value = "0.582"
$ws.Range("E23").Value = "  +1.40%  "
$ws.Range("D24").Value = "3.722.15"
$ws.Range("E24").Value = "  +0.42%  "
$ws.Range("D25").Value = "74.19"
$ws.Range("E25").Value = "  -0.30%  "
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("D27").Value = "0.0000115"
$ws.Range("E27").Value = "  -2.17%  "
$ws.Range("D28").Value = "8.07"
$ws.Range("E28").Value = "  +1.82%  "
$ws.Range("D29").Value = "1.65"
$ws.Range("E29").Value = "  +29.65%  "
$ws.Range("D30").Value = "8.66"
$ws.Range("E30").Value = "  +3.60%  "
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("E32").Value = "  +0.88%  "
$ws.Range("D33").Value = "3.584.57"
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("D34").Value = "24.31"
$ws.Range("E34").Value = "  +1.55%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  -0.35%  "
$ws.Range("D37").Value = "171.87"
$ws.Range("E37").Value = "  +2.31%  "
$ws.Range("D38").Value = "7.00"
$ws.Range("E38").Value = "  -1.56%  "
$ws.Range("D39").Value = "5.18"
$ws.Range("E39").Value = "  +2.63%  "
$ws.Range("E40").Value = "  +0.55%  "
$ws.Range("D41").Value = "0.0827"
$ws.Range("E41").Value = "  +2.47%  "
$ws.Range("D42").Value = "0.831"
$ws.Range("E42").Value = "  -0.48%  "
$ws.Range("E43").Value = "  -1.41%  "
$ws.Range("B44").Value = "ONDO"
$ws.Range("C44").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D44").Value = "1.26"
$ws.Range("E44").Value = "  +4.71%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Value = "43.11"
$ws.Range("E45").Value = "  +0.23%  "
$ws.Range("E46").Value = "  -0.11%  "
$ws.Range("D47").Value = "4.48"
$ws.Range("E47").Value = "  +0.33%  "
$ws.Range("D48").Value = "1.68"
$ws.Range("E48").Value = "  -1.05%  "
$ws.Range("D49").Value = "6.98"
$ws.Range("E49").Value = "  +1.95%  "
$ws.Range("D50").Value = "2.452.58"
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("D51").Value = "0.0269"
$ws.Range("E51").Value = "  +1.33%  "
